# Revised Project 1 Margins and Plots.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Margins (Project 1 Pars)")

# New margin / CI values (columns D, F, G) and updated significance labels (column E)
$rowData = @{
    2  = @{ D = 3.3075009999999998; F = 3.2543859999999998; G = 3.360617 }
    3  = @{ D = 3.4236140000000002; E = "NS"; F = 3.3363070000000001; G = 3.5109210000000002 }
    4  = @{ D = 3.3906540000000001; E = "NS"; F = 3.2583440000000001; G = 3.5229629999999998 }
    5  = @{ D = 3.4118170000000001; F = 3.3758189999999999; G = 3.447816 }
    6  = @{ D = 3.4547780000000001; E = "NS"; F = 3.3970009999999999; G = 3.5125549999999999 }
    7  = @{ D = 3.47167;            E = "NS"; F = 3.389154;            G = 3.5541849999999999 }
    8  = @{ D = 3.516133;           F = 3.4890840000000001; G = 3.543183 }
    9  = @{ D = 3.4859420000000001; E = "NS"; F = 3.4460709999999999; G = 3.5258129999999999 }
    10 = @{ D = 3.5526849999999999; E = "NS"; F = 3.493458;            G = 3.6119129999999999 }
    11 = @{ D = 3.6204489999999998; F = 3.586856;             G = 3.6540430000000002 }
    12 = @{ D = 3.5171060000000001; E = "NS"; F = 3.4683519999999999; G = 3.5658599999999998 }
    13 = @{ D = 3.6337009999999998; E = "NS"; F = 3.5465059999999999; G = 3.7208969999999999 }
    14 = @{ D = 3.7247650000000001; F = 3.6749010000000002; G = 3.7746300000000002 }
    15 = @{ D = 3.54827;            E = "NS"; F = 3.4727709999999998; G = 3.6237689999999998 }
    16 = @{ D = 3.7147169999999998; E = "NS"; F = 3.5765349999999998; G = 3.8528989999999999 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$r").Value = $vals.E
    }
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
}

# Update the selected range to reflect the author's final selection
$ws.Range("E15:E16").Select()
